$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 8 timestamp with refined float precision
$ws.Range("A8").Value = 45865.33355730324

# Append new row 9 with the latest sensor reading
$ws.Range("A9").Value = 45865.37525619288
$ws.Range("B9").Value = 2025
$ws.Range("C9").Value = 30
$ws.Range("D9").Value = 14.39
$ws.Range("E9").Value = 86.68000000000001
$ws.Range("F9").Value = 170.51
$ws.Range("G9").Value = 4.06
$ws.Range("H9").Value = "ESE"
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "09:00:22"

# Match formatting of row 8: A column uses style 2 (date/time number format)
$ws.Range("A9").NumberFormat = $ws.Range("A8").NumberFormat
